$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Milan"
$ws.Range("B3").Value = "Bern"
$ws.Range("A4").Value = "London"
$ws.Range("B4").Value = "Brisbane"

$ws.Range("B4").Select()
